$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing existing rows (and the
# hyperlink on the old A3) down by one.
$ws.Rows(3).Insert()

# New row 3 gets a whitespace-only value (four spaces).
$ws.Range("A3").Value = "    "

# The hyperlink that used to live on A3 now belongs on A4 (it moved
# down with the row insert, but the engine does not auto-relocate
# hyperlink anchors on row insert, so do it explicitly).
$ws.Range("A3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:#@%^%#$@#$@#.com")
$ws.Range("A4").Style = "Hyperlink"

# Match the active selection recorded in the saved workbook.
[void]$ws.Range("K8").Select()
